# Re-order several match rows that share the same kickoff timestamp
# (columns F:V move as a whole block between rows; columns A:E -- index,
# country, competition, season, kickoff datetime -- stay put on their own
# row) and append one brand-new match row (152) at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> row whose F:V content it should receive
$map = @{
    75  = 76
    76  = 75
    82  = 84
    83  = 85
    84  = 82
    85  = 83
    91  = 92
    92  = 93
    93  = 91
    94  = 95
    95  = 94
    103 = 104
    104 = 105
    105 = 106
    106 = 103
    126 = 128
    127 = 126
    128 = 127
    133 = 134
    134 = 133
    140 = 141
    141 = 140
    144 = 145
    145 = 144
    150 = 151
    151 = 150
}

# Snapshot every F:V block referenced above BEFORE any writes happen,
# so rows that feed more than one destination are not clobbered mid-way.
$snapshot = @{}
foreach ($r in $map.Keys) {
    $snapshot[$r] = $ws.Range("F$r`:V$r").Value2
}

foreach ($r in $map.Keys) {
    $src = $map[$r]
    $ws.Range("F$r`:V$r").Value2 = $snapshot[$src]
}

# Append the new match row 152.
$ws.Range("A152").Value2 = 151
$ws.Range("B152").Value2 = "italy"
$ws.Range("C152").Value2 = "serie-c-group-c"
$ws.Range("D152").Value2 = "2023-2024"
$ws.Range("E152").Value2 = 45262.67708333334
$ws.Range("F152").Value2 = "Monterosi"
$ws.Range("G152").Value2 = 0
$ws.Range("H152").Value2 = "ACR Messina"
$ws.Range("I152").Value2 = 2
$ws.Range("J152").Value2 = 2.24
$ws.Range("K152").Value2 = "30/11/2023 09:13"
$ws.Range("L152").Value2 = 2.29
$ws.Range("M152").Value2 = "02/12/2023 14:18"
$ws.Range("N152").Value2 = 2.92
$ws.Range("O152").Value2 = "30/11/2023 09:13"
$ws.Range("P152").Value2 = 3.08
$ws.Range("Q152").Value2 = "02/12/2023 14:18"
$ws.Range("R152").Value2 = 3.28
$ws.Range("S152").Value2 = "30/11/2023 09:13"
$ws.Range("T152").Value2 = 3.36
$ws.Range("U152").Value2 = "02/12/2023 14:18"
$ws.Range("V152").Value2 = "https://www.betexplorer.com/football/italy/serie-c-group-c/monterosi-acr-messina/AHS2udnC/"

# Match the styles used by the other data rows (bold/boxed index column,
# date-formatted kickoff column) by copying formats from the row above.
$ws.Range("A151").Copy() | Out-Null
$ws.Range("A152").PasteSpecial(-4122) | Out-Null

$ws.Range("E151").Copy() | Out-Null
$ws.Range("E152").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
